$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 8).Value = 1318.8085
$ws.Cells.Item(17, 9).Value = 1190
$ws.Cells.Item(17, 10).Value = 1349.3158
$ws.Cells.Item(17, 11).Value = 3570
$ws.Cells.Item(17, 12).Value = 4047.9474
$ws.Cells.Item(17, 13).Value = -3402
$ws.Cells.Item(17, 14).Value = -4383.9474
$ws.Cells.Item(58, 8).Value = 2837.7778
$ws.Cells.Item(58, 10).Value = 3898
$ws.Cells.Item(58, 12).Value = 11694
$ws.Cells.Item(58, 14).Value = -11994
$ws.Cells.Item(74, 8).Value = 3569
$ws.Cells.Item(74, 9).Value = 2649.7144
$ws.Cells.Item(74, 11).Value = 2649.7144
$ws.Cells.Item(74, 13).Value = -1713.7144
$ws.Cells.Item(77, 8).Value = 3569
$ws.Cells.Item(77, 9).Value = 2649.7144
$ws.Cells.Item(77, 11).Value = 13248.572
$ws.Cells.Item(77, 13).Value = -8568.572
$ws.Cells.Item(123, 8).Value = 145555
$ws.Cells.Item(123, 10).Value = 145555
$ws.Cells.Item(123, 12).Value = 145555
$ws.Cells.Item(123, 14).Value = -155355
$ws.Cells.Item(126, 8).Value = 145555
$ws.Cells.Item(126, 10).Value = 145555
$ws.Cells.Item(126, 12).Value = 145555
$ws.Cells.Item(126, 14).Value = -155435
$ws.Cells.Item(129, 8).Value = 3013.5715
$ws.Cells.Item(129, 10).Value = 5999
$ws.Cells.Item(129, 12).Value = 17997
$ws.Cells.Item(129, 14).Value = -27997
$ws.Cells.Item(137, 8).Value = 913990.6
$ws.Cells.Item(137, 9).Value = 1216949.8
$ws.Cells.Item(137, 10).Value = 5113
$ws.Cells.Item(137, 11).Value = 3650849.4
$ws.Cells.Item(137, 12).Value = 15339
$ws.Cells.Item(137, 13).Value = -3648299.4
$ws.Cells.Item(137, 14).Value = -20439
$ws.Cells.Item(138, 8).Value = 161606.58
$ws.Cells.Item(138, 10).Value = 5689.3184
$ws.Cells.Item(138, 12).Value = 17067.9552
$ws.Cells.Item(138, 14).Value = -27347.9552

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 13038.018
$ws.Cells.Item(32, 9).Value = 11262.115
$ws.Cells.Item(32, 11).Value = 11262.115
$ws.Cells.Item(32, 13).Value = -10975.115
$ws.Cells.Item(48, 8).Value = 227244
$ws.Cells.Item(48, 10).Value = 227244
$ws.Cells.Item(48, 12).Value = 227244
$ws.Cells.Item(48, 14).Value = -228012
$ws.Cells.Item(51, 8).Value = 50000
$ws.Cells.Item(51, 10).Value = 50000
$ws.Cells.Item(51, 12).Value = 50000
$ws.Cells.Item(51, 14).Value = -51512
$ws.Cells.Item(74, 8).Value = 3634.8833
$ws.Cells.Item(74, 9).Value = 10753.75
$ws.Cells.Item(74, 10).Value = 1855.1666
$ws.Cells.Item(74, 11).Value = 10753.75
$ws.Cells.Item(74, 12).Value = 1855.1666
$ws.Cells.Item(74, 13).Value = -9879.75
$ws.Cells.Item(74, 14).Value = -3603.1666
$ws.Cells.Item(77, 8).Value = 3634.8833
$ws.Cells.Item(77, 9).Value = 10753.75
$ws.Cells.Item(77, 10).Value = 1855.1666
$ws.Cells.Item(77, 11).Value = 53768.75
$ws.Cells.Item(77, 12).Value = 9275.833000000001
$ws.Cells.Item(77, 13).Value = -49400.75
$ws.Cells.Item(77, 14).Value = -18011.833
$ws.Cells.Item(109, 8).Value = 100000
$ws.Cells.Item(109, 10).Value = 100000
$ws.Cells.Item(109, 12).Value = 100000
$ws.Cells.Item(109, 14).Value = -102774
$ws.Cells.Item(132, 8).Value = 2047.76
$ws.Cells.Item(132, 9).Value = 1934.75
$ws.Cells.Item(132, 10).Value = 2499.8
$ws.Cells.Item(132, 11).Value = 5804.25
$ws.Cells.Item(132, 12).Value = 7499.400000000001
$ws.Cells.Item(132, 13).Value = -3274.25
$ws.Cells.Item(132, 14).Value = -12559.4
$ws.Cells.Item(133, 8).Value = 79998.664
$ws.Cells.Item(133, 10).Value = 79998.664
$ws.Cells.Item(133, 12).Value = 79998.664
$ws.Cells.Item(133, 14).Value = -85058.664

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 3655.75
$ws.Cells.Item(20, 9).Value = 3055.1538
$ws.Cells.Item(20, 10).Value = 4771.143
$ws.Cells.Item(20, 11).Value = 3055.1538
$ws.Cells.Item(20, 12).Value = 4771.143
$ws.Cells.Item(20, 13).Value = -2808.1538
$ws.Cells.Item(20, 14).Value = -5265.143
$ws.Cells.Item(86, 8).Value = 7976.846
$ws.Cells.Item(86, 9).Value = 8519.272000000001
$ws.Cells.Item(86, 11).Value = 8519.272000000001
$ws.Cells.Item(86, 13).Value = -7396.272000000001
$ws.Cells.Item(89, 8).Value = 7976.846
$ws.Cells.Item(89, 9).Value = 8519.272000000001
$ws.Cells.Item(89, 11).Value = 42596.36
$ws.Cells.Item(89, 13).Value = -36980.36
$ws.Cells.Item(122, 8).Value = 77780
$ws.Cells.Item(122, 10).Value = 77780
$ws.Cells.Item(122, 12).Value = 77780
$ws.Cells.Item(122, 14).Value = -87580

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(22, 8).Value = 406.66666
$ws.Cells.Item(22, 10).Value = 537.3333
$ws.Cells.Item(22, 12).Value = 537.3333
$ws.Cells.Item(22, 14).Value = -1237.3333
$ws.Cells.Item(31, 8).Value = 3758.6365
$ws.Cells.Item(31, 9).Value = 2193.2354
$ws.Cells.Item(31, 11).Value = 2193.2354
$ws.Cells.Item(31, 13).Value = -1898.2354
$ws.Cells.Item(34, 8).Value = 3758.6365
$ws.Cells.Item(34, 9).Value = 2193.2354
$ws.Cells.Item(34, 11).Value = 2193.2354
$ws.Cells.Item(34, 13).Value = -1991.2354
$ws.Cells.Item(141, 8).Value = 308648.78
$ws.Cells.Item(141, 10).Value = 424441.16
$ws.Cells.Item(141, 12).Value = 424441.16
$ws.Cells.Item(141, 14).Value = -434801.16

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 8).Value = 118.72727
$ws.Cells.Item(2, 9).Value = 22.285715
$ws.Cells.Item(2, 10).Value = 287.5
$ws.Cells.Item(2, 11).Value = 133.71429
$ws.Cells.Item(2, 12).Value = 1725
$ws.Cells.Item(2, 13).Value = -20.71429000000001
$ws.Cells.Item(2, 14).Value = -1951
$ws.Cells.Item(80, 8).Value = 166934530
$ws.Cells.Item(80, 9).Value = 500000000
$ws.Cells.Item(80, 10).Value = 401800
$ws.Cells.Item(80, 11).Value = 1500000000
$ws.Cells.Item(80, 12).Value = 1205400
$ws.Cells.Item(80, 13).Value = -1499999064
$ws.Cells.Item(80, 14).Value = -1207272
$ws.Cells.Item(81, 8).Value = 201779.6
$ws.Cells.Item(81, 9).Value = 500249.5
$ws.Cells.Item(81, 11).Value = 1500748.5
$ws.Cells.Item(81, 13).Value = -1499625.5
$ws.Cells.Item(83, 8).Value = 166934530
$ws.Cells.Item(83, 9).Value = 500000000
$ws.Cells.Item(83, 10).Value = 401800
$ws.Cells.Item(83, 11).Value = 4500000000
$ws.Cells.Item(83, 12).Value = 3616200
$ws.Cells.Item(83, 13).Value = -4499995320
$ws.Cells.Item(83, 14).Value = -3625560
$ws.Cells.Item(84, 8).Value = 201779.6
$ws.Cells.Item(84, 9).Value = 500249.5
$ws.Cells.Item(84, 11).Value = 4502245.5
$ws.Cells.Item(84, 13).Value = -4496629.5
$ws.Cells.Item(130, 8).Value = 14007.111
$ws.Cells.Item(130, 9).Value = 2021.3334
$ws.Cells.Item(130, 11).Value = 6064.0002
$ws.Cells.Item(130, 13).Value = -1044.0002
$ws.Cells.Item(137, 8).Value = 11260.357
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 11260.357
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).ClearContents()
$ws.Cells.Item(137, 13).Value = 33781.071
$ws.Cells.Item(137, 14).Value = -43981.071

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(45, 8).Value = 33333.332
$ws.Cells.Item(45, 10).Value = 33333.332
$ws.Cells.Item(45, 12).Value = 33333.332
$ws.Cells.Item(45, 14).Value = -34451.332
$ws.Cells.Item(101, 8).Value = 35000
$ws.Cells.Item(101, 10).Value = 35000
$ws.Cells.Item(101, 12).Value = 35000
$ws.Cells.Item(101, 14).Value = -41490
$ws.Cells.Item(122, 8).Value = 23611
$ws.Cells.Item(122, 9).Value = 18749.834
$ws.Cells.Item(122, 11).Value = 56249.50199999999
$ws.Cells.Item(122, 13).Value = -53799.50199999999
$ws.Cells.Item(132, 8).Value = 3838.4194
$ws.Cells.Item(132, 9).Value = 3635.68
$ws.Cells.Item(132, 10).Value = 4683.1665
$ws.Cells.Item(132, 11).Value = 10907.04
$ws.Cells.Item(132, 12).Value = 14049.4995
$ws.Cells.Item(132, 13).Value = -8377.039999999999
$ws.Cells.Item(132, 14).Value = -19109.4995

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(82, 8).Value = 2114.111
$ws.Cells.Item(82, 10).Value = 2204.5
$ws.Cells.Item(82, 12).Value = 2204.5
$ws.Cells.Item(82, 14).Value = -2926.5
$ws.Cells.Item(85, 8).Value = 2114.111
$ws.Cells.Item(85, 10).Value = 2204.5
$ws.Cells.Item(85, 12).Value = 2204.5
$ws.Cells.Item(85, 14).Value = -4700.5
$ws.Cells.Item(93, 8).Value = 4617.6665
$ws.Cells.Item(93, 9).Value = 4693.6924
$ws.Cells.Item(93, 11).Value = 4693.6924
$ws.Cells.Item(93, 13).Value = -3445.6924
$ws.Cells.Item(101, 8).Value = 32501
$ws.Cells.Item(101, 10).Value = 32501
$ws.Cells.Item(101, 12).Value = 32501
$ws.Cells.Item(101, 14).Value = -38991
$ws.Cells.Item(106, 8).Value = 13500
$ws.Cells.Item(106, 10).Value = 13500
$ws.Cells.Item(106, 12).Value = 13500
$ws.Cells.Item(106, 14).Value = -16024
$ws.Cells.Item(122, 8).Value = 7731.8887
$ws.Cells.Item(122, 9).Value = 9417.714
$ws.Cells.Item(122, 10).Value = 6659.091
$ws.Cells.Item(122, 11).Value = 28253.142
$ws.Cells.Item(122, 12).Value = 19977.273
$ws.Cells.Item(122, 13).Value = -25803.142
$ws.Cells.Item(122, 14).Value = -24877.273
$ws.Cells.Item(132, 8).Value = 351022.22
$ws.Cells.Item(132, 10).Value = 4186.1113
$ws.Cells.Item(132, 12).Value = 12558.3339
$ws.Cells.Item(132, 14).Value = -17618.3339

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(103, 8).Value = 600
$ws.Cells.Item(103, 9).Value = 600
$ws.Cells.Item(103, 11).Value = 600
$ws.Cells.Item(103, 13).Value = 572
$ws.Cells.Item(125, 8).Value = 120000
$ws.Cells.Item(125, 10).Value = 120000
$ws.Cells.Item(125, 12).Value = 120000
$ws.Cells.Item(125, 14).Value = -129840
$ws.Cells.Item(126, 8).Value = 37530
$ws.Cells.Item(126, 9).Value = 46707.223
$ws.Cells.Item(126, 10).Value = 9998.333000000001
$ws.Cells.Item(126, 11).Value = 140121.669
$ws.Cells.Item(126, 12).Value = 29994.999
$ws.Cells.Item(126, 13).Value = -137651.669
$ws.Cells.Item(126, 14).Value = -34934.999
$ws.Cells.Item(132, 8).Value = 11380.726
$ws.Cells.Item(132, 9).Value = 13944.857
$ws.Cells.Item(132, 11).Value = 41834.571
$ws.Cells.Item(132, 13).Value = -39304.571
$ws.Cells.Item(133, 8).Value = 100000
$ws.Cells.Item(133, 10).Value = 100000
$ws.Cells.Item(133, 12).Value = 100000
$ws.Cells.Item(133, 14).Value = -110120

Write-Output "Applied all Siren_Profits updates"